$d = $word.ActiveDocument

# --- Delete paragraphs that no longer exist in the target (highest index first) ---
# 7: "Capital geht z.T. ins Minus, ..."
$d.Paragraphs(7).Range.Delete()
# 6: "Verkaufspreis muss geändert werden"
$d.Paragraphs(6).Range.Delete()
# 2: "Pflanzen wachsen lassen"
$d.Paragraphs(2).Range.Delete()

# --- Paragraph 1: "State-Bildchen in State" -> extended English-tagged text ---
$p1 = $d.Paragraphs(1).Range
$p1.Text = "State-Bildchen in State (function, interval)?"
$p1.LanguageID = "en-US"

# --- Paragraph (now index 5): "neues Klassendiagramm" -> append " + AD" ---
$p5 = $d.Paragraphs(5).Range
$p5.Text = "neues Klassendiagramm + AD"
